# Generate Report for Handback
#
# - "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it
#   appears (Overview!E2:F2/E3:F3, zh-cn!C2:C3, de-de!C2:C3).
# - zh-cn / de-de sheets: fill in the (previously empty) "Latest Target
#   File" (I) / "Latest Handback File" (J) columns for rows 2-3, turn I2/I3
#   into a hyperlink to a.md (mirroring the existing A2 hyperlink), and
#   bump the "Latest Handback DateTime" (K) timestamp.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status columns E/F for rows 2 and 3.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $overview.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# ---------------------------------------------------------------------
# Per-language handback sheets.
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; XlfFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; HandbackDateTime = "2016-11-14 06:16:36" },
    @{ Name = "de-de"; XlfFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; HandbackDateTime = "2016-11-14 06:16:54" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (C) for rows 2 and 3.
    foreach ($addr in @("C2", "C3")) {
        $cell = $ws.Range($addr)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }

    # Address used by the existing "a.md" hyperlink in column A - reuse it
    # for the new column-I hyperlinks.
    $aLink = $ws.Hyperlinks | Where-Object { $_.Range.Address() -eq '$A$2' } | Select-Object -First 1
    $aMdAddress = $aLink.Address

    foreach ($row in @(2, 3)) {
        $iCell = $ws.Range("I$row")
        $jCell = $ws.Range("J$row")
        $kCell = $ws.Range("K$row")

        # Latest Target File -> "a.md" hyperlink (mirrors column A).
        $ws.Hyperlinks.Add($iCell, $aMdAddress, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
        $iCell.Style = "HyperLink"

        # Latest Handback File.
        $jCell.Value = $lang.XlfFile

        # Latest Handback DateTime.
        $kCell.Value = $lang.HandbackDateTime
    }
}

Write-Host "Handback report generated."
